# sprint_backlog.xlsx - "andrat i basicview funktioner"
# Update the week 48-49 actuals (rows 60, 64-67) to DONE with end dates / comments,
# and append the week 50-51 block (rows 69-78) following the same weekly template
# used for 46-47 (rows 48-55) and 48-49 (rows 57-67).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 60: "Find more relevant articles about HARQ and LA" -> DONE
# ---------------------------------------------------------------------------
$ws.Range("D62").Copy() | Out-Null
$ws.Range("D60").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D60").Value = "DONE"

$ws.Range("G64").Copy() | Out-Null
$ws.Range("G60").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G60").Value = 41975

$ws.Range("H60").Value = "The articles is now about BLER vs Throughput"

# ---------------------------------------------------------------------------
# Row 64: "Demo for Ola and Johannes" -> DONE, real enddate
# ---------------------------------------------------------------------------
$ws.Range("D64").Copy() | Out-Null
$ws.Range("D64").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D62").Copy() | Out-Null
$ws.Range("D64").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D64").Value = "DONE"
$ws.Range("G64").Value = 41975

# ---------------------------------------------------------------------------
# Row 65: "plan the new analysation BLER target" -> DONE, add enddate
# ---------------------------------------------------------------------------
$ws.Range("D62").Copy() | Out-Null
$ws.Range("D65").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D65").Value = "DONE"

$ws.Range("G64").Copy() | Out-Null
$ws.Range("G65").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G65").Value = 41975

# ---------------------------------------------------------------------------
# Row 66: "Formulate the question ..." -> DONE, real enddate
# ---------------------------------------------------------------------------
$ws.Range("D62").Copy() | Out-Null
$ws.Range("D66").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D66").Value = "DONE"
$ws.Range("G66").Value = 41977

# ---------------------------------------------------------------------------
# Row 67: "find a scientific base (articles) ..." -> DONE, add enddate
# ---------------------------------------------------------------------------
$ws.Range("D62").Copy() | Out-Null
$ws.Range("D67").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D67").Value = "DONE"
$ws.Range("G67").Value = 41975

# ---------------------------------------------------------------------------
# Week 50-51 block: rows 69-78 (mirrors the week 48-49 block at rows 57-67)
# ---------------------------------------------------------------------------

# Row 69 - new week header/divider row (clone of row 57's template)
$ws.Range("B57:H57").Copy() | Out-Null
$ws.Range("B69:H69").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B69").Value = "50-51"
$ws.Range("C69").Value = "create an realtime analyser in the plugin"
$ws.Range("D69").Value = "TODO?"
$ws.Range("E69").Value = "Razmus"

# Row 70 - carried over ONGOING task (clone of row 58)
$ws.Range("B58:H58").Copy() | Out-Null
$ws.Range("B70:H70").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C70").Value = "Writing the end report."
$ws.Range("D70").Value = "ONGOING"
$ws.Range("E70").Value = "Paul & Razmus"
$ws.Range("F70").Value = 41902

# Row 71 - carried over ONGOING task (clone of row 59)
$ws.Range("B59:H59").Copy() | Out-Null
$ws.Range("B71:H71").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C71").Value = "build the analyzer-plugin to logtool in Java."
$ws.Range("D71").Value = "ONGOING"
$ws.Range("E71").Value = "Razmus"
$ws.Range("F71").Value = 41911

# Row 72 - carried over ONGOING task (clone of row 61)
$ws.Range("B61:H61").Copy() | Out-Null
$ws.Range("B72:H72").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C72").Value = "fix the optional tab in the analyser"
$ws.Range("D72").Value = "ONGOING"
$ws.Range("E72").Value = "Razmus"
$ws.Range("F72").Value = 41953

# Row 73 - new TODO task, demo date, comment (no names)
$ws.Range("B65:H65").Copy() | Out-Null
$ws.Range("B73:H73").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B73").ClearContents()
$ws.Range("C73").Value = "Have a demo for the tester at IODT"
$ws.Range("D69").Copy() | Out-Null
$ws.Range("D73").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D73").Value = "TODO"
$ws.Range("E73").ClearContents()
$ws.Range("F73").Value = 41981
$ws.Range("G73").ClearContents()
$ws.Range("H73").Value = "The Demo is between 14-15 at Monday"

# Row 74 - new STARTED task
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C74").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C74").Value = "add automatic analysation for the tool"
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D74").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D74").Value = "STARTED"
$ws.Range("E74").Value = "Paul"

# Row 75 - new TODO task with comment
$ws.Range("C75").Style = $ws.Range("C5").Style
$ws.Range("C75").Value = "Create a servey for the testers"
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D75").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D75").Value = "TODO"
$ws.Range("E75").Value = "Razmus"
$ws.Range("H75").Value = "The servey will be question on how good the tool is, easy to use etc..."

# Row 76 - new TODO task with a preliminary (text) date in F
$ws.Range("C76").Style = $ws.Range("C5").Style
$ws.Range("C76").Value = "Draw traces with different BLER settings."
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D76").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D76").Value = "TODO"
$ws.Range("E76").Value = "Paul"
$ws.Range("F76").Value = "prel 5-dec"

# Row 77 - new TODO task with comment
$ws.Range("C77").Style = $ws.Range("C5").Style
$ws.Range("C77").Value = "Do analysation on the trace data + test the automatic analyse tool"
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D77").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D77").Value = "TODO"
$ws.Range("E77").Value = "Paul"
$ws.Range("H77").Value = "look at which the optimal bler target is + calculate it with the automation function"

# Row 78 - new STARTED task with start date and comment
$ws.Range("C78").Style = $ws.Range("C7").Style
$ws.Range("C78").Value = "Bug fixing + redesign of the tool"
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D78").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D78").Value = "STARTED"
$ws.Range("E78").Value = "(Paul) & Razmus"
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F78").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F78").Value = 41973
$ws.Range("H78").Value = "Fix bugs, optimize code, change view layout etc."

# ---------------------------------------------------------------------------
# Cosmetic sheet-level tweaks
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 15.585
$ws.Range("C80").Select() | Out-Null
